# Add a "Kỹ năng" (Skill) column to the student import template, between
# the existing "Ngày ra trường" (H) and "Mô tả" (I -> J) columns, and
# populate the sample row with "php, laravel". This shifts the former
# column I ("Mô tả") to column J.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at I; this pushes the old column I ("Mô tả") to J
# and makes the new column I inherit formatting from its left neighbor H
# (yellow header fill / date-ish column styles), matching Excel's default
# "Insert" behaviour.
$ws.Columns("I:I").Insert() | Out-Null

# Match the new column's width to column H's width (both render as 17 in
# the saved column width units).
$ws.Columns("I:I").ColumnWidth = 16.17

# New header + sample value for the inserted "Kỹ năng" column.
$ws.Range("I1").Value = "Kỹ năng"
$ws.Range("I2").Value = "php, laravel"

# Update the remembered selection to the cell that shifted from H7 to J7.
$ws.Range("J7").Select() | Out-Null
